$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.441.06"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.573.29"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.74"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3730"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.88"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3391"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.142"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.34"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.005"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.954"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "1.577.05"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.90"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06747"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.295"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").Value = "22.427.53"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.349"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.671"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.08"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.28"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.000"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.38"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "1.751.86"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.175"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.976"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.824"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08385"
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02477"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2282"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06525"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.459"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.27"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.94"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.815"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5806"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.56"
$ws.Range("E48").Value = "  +3.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.075"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").Value = "  -7.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07327"
$ws.Range("E51").Value = "  +0.00%  "
